$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 54 data, matching the style of the row above it (row 53)
$ws.Range("A53").Copy($ws.Range("A54"))
$ws.Range("A54").Value = 45986

$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = -2.06674933094535
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = -0.9969640812590996
